# Applies the "Add files via upload" update to "fii dii.xlsx":
# appends two new daily blocks (11/06/2024-derived data ending with dates
# "12/06/2024" and "13/06/2024") to Sheet1, rows 839-919, columns A-J.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('A839').Value = 'Buying Opportunity'
$ws.Range('B839').Value = 'support Zone'
$ws.Range('C839').Value = 'long buildup'
$ws.Range('D839').Value = 'Short buildup'
$ws.Range('E839').Value = 'FII ENTERING'

$ws.Range('A840').Value = 'ABSLAMC'
$ws.Range('B840').Value = 'ARE&M'
$ws.Range('C840').Value = 'IEX'
$ws.Range('E840').Value = 'LICHSGFIN'
$ws.Range('F840').Value = 639.45
$ws.Range('G840').Value = 1342.45
$ws.Range('H840').Value = 173.61
$ws.Range('J840').Value = 715.4

$ws.Range('A841').Value = 'AIAENG'
$ws.Range('B841').Value = 'ATGL'
$ws.Range('C841').Value = 'LICHSGFIN'
$ws.Range('E841').Value = 'PEL'
$ws.Range('F841').Value = 3825.85
$ws.Range('G841').Value = 951.7
$ws.Range('H841').Value = 715.4
$ws.Range('J841').Value = 865.6

$ws.Range('A842').Value = 'ALPHAGEO'
$ws.Range('B842').Value = 'FINPIPE'
$ws.Range('F842').Value = 389.4
$ws.Range('G842').Value = 339.25

$ws.Range('A843').Value = 'AMBER'
$ws.Range('B843').Value = 'IRMENERGY'
$ws.Range('F843').Value = 4077.4
$ws.Range('G843').Value = 464.95

$ws.Range('A844').Value = 'ATUL'
$ws.Range('B844').Value = 'METROBRAND'
$ws.Range('F844').Value = 6192.4
$ws.Range('G844').Value = 1124.55

$ws.Range('A845').Value = 'BAJAJELEC'
$ws.Range('F845').Value = 1017.15

$ws.Range('A846').Value = 'CARBORUNIV'
$ws.Range('F846').Value = 1705.7

$ws.Range('A847').Value = 'CHEMFAB'
$ws.Range('F847').Value = 688.1

$ws.Range('A848').Value = 'CLEDUCATE'
$ws.Range('F848').Value = 88.43000000000001

$ws.Range('A849').Value = 'CONCOR'
$ws.Range('F849').Value = 1141.35

$ws.Range('A850').Value = 'DBOL'
$ws.Range('F850').Value = 156.75

$ws.Range('A851').Value = 'DCMSRIND'
$ws.Range('F851').Value = 221.33

$ws.Range('A852').Value = 'DEEPAKNTR'
$ws.Range('F852').Value = 2351.8

$ws.Range('A853').Value = 'DREDGECORP'
$ws.Range('F853').Value = 1129.2

$ws.Range('A854').Value = 'FLUOROCHEM'
$ws.Range('F854').Value = 3222.85

$ws.Range('A855').Value = 'GNA'
$ws.Range('F855').Value = 405.3

$ws.Range('A856').Value = 'GRAVITA'
$ws.Range('F856').Value = 1332.85

$ws.Range('A857').Value = 'HATHWAY'
$ws.Range('F857').Value = 22.26

$ws.Range('A858').Value = 'IBULHSGFIN'
$ws.Range('F858').Value = 170.59

$ws.Range('A859').Value = 'IGARASHI'
$ws.Range('F859').Value = 515.7

$ws.Range('A860').Value = 'INDBANK'
$ws.Range('F860').Value = 49.77

$ws.Range('A861').Value = 'INDIANHUME'
$ws.Range('F861').Value = 352.3

$ws.Range('A862').Value = 'INDNIPPON'
$ws.Range('F862').Value = 713.05

$ws.Range('A863').Value = 'INFOMEDIA'
$ws.Range('F863').Value = 6.58

$ws.Range('A864').Value = 'JWL'
$ws.Range('F864').Value = 694.9

$ws.Range('A865').Value = 'KRITINUT'
$ws.Range('F865').Value = 113.62

$ws.Range('A866').Value = 'KRSNAA'
$ws.Range('F866').Value = 637.5

$ws.Range('A867').Value = 'LAMBODHARA'
$ws.Range('F867').Value = 151.69

$ws.Range('A868').Value = 'LICHSGFIN'
$ws.Range('F868').Value = 715.4

$ws.Range('A869').Value = 'MARINE'
$ws.Range('F869').Value = 132.17

$ws.Range('A870').Value = 'MCX'
$ws.Range('F870').Value = 3793.15

$ws.Range('A871').Value = 'MMTC'
$ws.Range('F871').Value = 77.97

$ws.Range('A872').Value = 'NAHARINDUS'
$ws.Range('F872').Value = 131.92

$ws.Range('A873').Value = 'NECLIFE'
$ws.Range('F873').Value = 33.7

$ws.Range('A874').Value = 'NINSYS'
$ws.Range('F874').Value = 520.7

$ws.Range('A875').Value = 'NURECA'
$ws.Range('F875').Value = 280.3

$ws.Range('A876').Value = 'ORIENTLTD'
$ws.Range('F876').Value = 114.21

$ws.Range('A877').Value = 'ORISSAMINE'
$ws.Range('F877').Value = 6760.7

$ws.Range('A878').Value = 'PATELENG'
$ws.Range('F878').Value = 66.65000000000001

$ws.Range('A879').Value = 'PDSL'
$ws.Range('F879').Value = 506.95

$ws.Range('A880').Value = 'PEL'
$ws.Range('F880').Value = 865.6

$ws.Range('A881').Value = 'PRAJIND'
$ws.Range('F881').Value = 634.7

$ws.Range('A882').Value = 'RADIANTCMS'
$ws.Range('F882').Value = 80.29000000000001

$ws.Range('A883').Value = 'RPOWER'
$ws.Range('F883').Value = 31.54

$ws.Range('A884').Value = 'SADHNANIQ'
$ws.Range('F884').Value = 64.91

$ws.Range('A885').Value = "'12/06/2024"

$ws.Range('A886').Value = 'Buying Opportunity'
$ws.Range('B886').Value = 'support Zone'
$ws.Range('C886').Value = 'long buildup'
$ws.Range('D886').Value = 'Short buildup'
$ws.Range('E886').Value = 'FII ENTERING'

$ws.Range('A887').Value = 'AGSTRA'
$ws.Range('B887').Value = 'ARE&M'
$ws.Range('E887').Value = 'OFSS'
$ws.Range('F887').Value = 80.37
$ws.Range('G887').Value = 1340.65
$ws.Range('J887').Value = 9658

$ws.Range('A888').Value = 'ALPHAGEO'
$ws.Range('B888').Value = 'ASHIMASYN'
$ws.Range('F888').Value = 412.3
$ws.Range('G888').Value = 35.73

$ws.Range('A889').Value = 'AMBIKCO'
$ws.Range('B889').Value = 'CAMPUS'
$ws.Range('F889').Value = 1645
$ws.Range('G889').Value = 284.25

$ws.Range('A890').Value = 'ASKAUTOLTD'
$ws.Range('B890').Value = 'COASTCORP'
$ws.Range('F890').Value = 342.3
$ws.Range('G890').Value = 232.72

$ws.Range('A891').Value = 'ASTRAMICRO'
$ws.Range('B891').Value = 'HMAAGRO'
$ws.Range('F891').Value = 926.5
$ws.Range('G891').Value = 54.68

$ws.Range('A892').Value = 'BCLIND'
$ws.Range('B892').Value = 'INTENTECH'
$ws.Range('F892').Value = 58.39
$ws.Range('G892').Value = 135.75

$ws.Range('A893').Value = 'BEL'
$ws.Range('B893').Value = 'JPOLYINVST'
$ws.Range('F893').Value = 300.9
$ws.Range('G893').Value = 761.9

$ws.Range('A894').Value = 'DCMSRIND'
$ws.Range('B894').Value = 'KARURVYSYA'
$ws.Range('F894').Value = 228.25
$ws.Range('G894').Value = 206.83

$ws.Range('A895').Value = 'DELTACORP'
$ws.Range('B895').Value = 'MANOMAY'
$ws.Range('F895').Value = 135.28
$ws.Range('G895').Value = 180.88

$ws.Range('A896').Value = 'DTIL'
$ws.Range('B896').Value = 'ROUTE'
$ws.Range('F896').Value = 208.27
$ws.Range('G896').Value = 1535.85

$ws.Range('A897').Value = 'FOSECOIND'
$ws.Range('F897').Value = 3906.65

$ws.Range('A898').Value = 'FUSION'
$ws.Range('F898').Value = 499.45

$ws.Range('A899').Value = 'GATEWAY'
$ws.Range('F899').Value = 104.48

$ws.Range('A900').Value = 'HGS'
$ws.Range('F900').Value = 854

$ws.Range('A901').Value = 'HINDCOMPOS'
$ws.Range('F901').Value = 457.4

$ws.Range('A902').Value = 'HONDAPOWER'
$ws.Range('F902').Value = 3125.2

$ws.Range('A903').Value = 'INDOSTAR'
$ws.Range('F903').Value = 254.91

$ws.Range('A904').Value = 'INFOMEDIA'
$ws.Range('F904').Value = 6.87

$ws.Range('A905').Value = 'JMFINANCIL'
$ws.Range('F905').Value = 84.37

$ws.Range('A906').Value = 'KSHITIJPOL'
$ws.Range('F906').Value = 7.03

$ws.Range('A907').Value = 'KUANTUM'
$ws.Range('F907').Value = 155.68

$ws.Range('A908').Value = 'LAMBODHARA'
$ws.Range('F908').Value = 157.12

$ws.Range('A909').Value = 'MARINE'
$ws.Range('F909').Value = 140.1

$ws.Range('A910').Value = 'MSTCLTD'
$ws.Range('F910').Value = 896.8

$ws.Range('A911').Value = 'NESCO'
$ws.Range('F911').Value = 981.5

$ws.Range('A912').Value = 'OFSS'
$ws.Range('F912').Value = 9658

$ws.Range('A913').Value = 'PARAS'
$ws.Range('F913').Value = 964.1

$ws.Range('A914').Value = 'PRAJIND'
$ws.Range('F914').Value = 685.7

$ws.Range('A915').Value = 'QUESS'
$ws.Range('F915').Value = 619.3

$ws.Range('A916').Value = 'RANASUG'
$ws.Range('F916').Value = 26.38

$ws.Range('A917').Value = 'RELINFRA'
$ws.Range('F917').Value = 203.15

$ws.Range('A918').Value = 'ROSSARI'
$ws.Range('F918').Value = 755.6

$ws.Range('A919').Value = "'13/06/2024"
